$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assigns a value while forcing text storage (prevents Excel from
# auto-converting numeric-looking strings, e.g. "7.00" or "0.500", into real
# numbers, which would silently drop significant trailing/leading zeros) and
# restores the cell's original style afterwards so no spurious style change
# (e.g. a "Number Stored as Text" quote-prefix style) is introduced.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.Value = "'" + $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = '66.015.50'
$ws.Range("E2").Value = '  +1.40%  '
$ws.Range("D3").Value = '3.593.42'
$ws.Range("E3").Value = '  +1.41%  '
Set-TextValue "D4" '0.999'
$ws.Range("E4").Value = '  -0.04%  '
Set-TextValue "D5" '602.27'
Set-TextValue "D6" '138.87'
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("D7").Value = '3.589.45'
$ws.Range("E7").Value = '  +1.33%  '
$ws.Range("E8").Value = '  -0.05%  '
Set-TextValue "D9" '0.500'
$ws.Range("E9").Value = '  +1.03%  '
$ws.Range("E10").Value = '  +1.36%  '
$ws.Range("E11").Value = '  +4.06%  '
Set-TextValue "D12" '0.392'
$ws.Range("E12").Value = '  +1.48%  '
$ws.Range("D13").Value = '4.203.45'
$ws.Range("E13").Value = '  +1.47%  '
Set-TextValue "D14" '28.18'
$ws.Range("E14").Value = '  +2.90%  '
$ws.Range("E15").Value = '  +1.73%  '
$ws.Range("D16").Value = '3.585.23'
$ws.Range("E16").Value = '  +1.10%  '
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '66.067.03'
$ws.Range("E18").Value = '  +1.49%  '
Set-TextValue "D19" '10.03'
$ws.Range("E19").Value = '  -0.36%  '
Set-TextValue "D20" '14.63'
$ws.Range("E20").Value = '  +2.51%  '
Set-TextValue "D21" '5.84'
$ws.Range("E21").Value = '  -1.04%  '
Set-TextValue "D22" '397.60'
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("E23").Value = '  +2.43%  '
$ws.Range("D24").Value = '3.734.62'
$ws.Range("E24").Value = '  +1.43%  '
Set-TextValue "D25" '74.97'
$ws.Range("E25").Value = '  +1.79%  '
$ws.Range("E26").Value = '  +0.10%  '
Set-TextValue "D27" '0.0000119'
$ws.Range("E27").Value = '  +4.22%  '
Set-TextValue "D28" '8.09'
$ws.Range("E28").Value = '  +3.26%  '
Set-TextValue "D29" '1.65'
$ws.Range("E29").Value = '  +20.89%  '
$ws.Range("E30").Value = '  +3.97%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue "D31" '1.01'
$ws.Range("E31").Value = '  +0.83%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D32" '2.33'
$ws.Range("E32").Value = '  +1.90%  '
$ws.Range("D33").Value = '3.594.58'
$ws.Range("E33").Value = '  +0.94%  '
Set-TextValue "D34" '24.57'
$ws.Range("E34").Value = '  +2.88%  '
Set-TextValue "D35" '0.149'
$ws.Range("E35").Value = '  +3.22%  '
$ws.Range("E37").Value = '  +7.70%  '
$ws.Range("E38").Value = '  +1.71%  '
Set-TextValue "D39" '7.00'
$ws.Range("E39").Value = '  +0.64%  '
Set-TextValue "D40" '169.63'
$ws.Range("E40").Value = '  -0.13%  '
Set-TextValue "D41" '0.0839'
$ws.Range("E41").Value = '  +4.43%  '
$ws.Range("E42").Value = '  +1.99%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D43" '26.27'
$ws.Range("E43").Value = '  -0.80%  '
$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue "D44" '1.27'
$ws.Range("E44").Value = '  +5.92%  '
Set-TextValue "D45" '43.09'
$ws.Range("E45").Value = '  +1.34%  '
Set-TextValue "D46" '4.54'
$ws.Range("E46").Value = '  +2.08%  '
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("E48").Value = '  +2.43%  '
$ws.Range("E49").Value = '  +1.56%  '
$ws.Range("D50").Value = '2.445.96'
$ws.Range("E50").Value = '  +2.55%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D51" '22.73'
$ws.Range("E51").Value = '  +7.47%  '

Write-Host "Applied cryptos.xlsx update on $(Get-Date)"
